$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 651, shifting existing rows 651:692 down to 652:693
$ws.Rows.Item(651).Insert()

# Populate the newly inserted row 651 with the new data point (2026/01/16 entry)
$ws.Cells.Item(651, 1).Value = "'2026/01/16"
$ws.Cells.Item(651, 1).Style = "Normal"
$ws.Cells.Item(651, 2).Value = "金"
$ws.Cells.Item(651, 3).Value = 19
$ws.Cells.Item(651, 4).Value = 201
